$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assessment of the IDFlakies detection tool performance:
# Recall OD dropped slightly (0.78 -> 0.77) and execution Time improved (3146 -> 2495)
$ws.Range("B5").Value = 0.77
$ws.Range("D5").Value = 2495

$ws.Range("B29").Value = 0.77
$ws.Range("D29").Value = 2495

# Re-enter these formula blocks as single range assignments so Excel stores
# them as shared formulas (matches how Excel re-saves after a fill/edit)
$ws.Range("B37:D39").Formula = "=POWER(B29,2)"
$ws.Range("B78:D80").Formula = "=B54"

# Update the view: scroll position / active selection moved
$ws.Range("D30").Select()
